$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 19
$ws.Range("I2").Value = 32
$ws.Range("J2").Value = 144
$ws.Range("L2").Value = 43
$ws.Range("M2").Value = 2
$ws.Range("N2").Value = 25
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 2
$ws.Range("R2").Value = 1
$ws.Range("S2").Value = 15
$ws.Range("T2").Value = 20
$ws.Range("U2").Value = 1
$ws.Range("V2").Value = 270
$ws.Range("Y2").Value = 2
